# Auto-generated: apply scraped scheduled-runner price/profit updates to Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 102.333336
$ws.Range("I4").Value = 102.333336
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 102.333336
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 11.666664
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 1093
$ws.Range("I11").Value = 1093
$ws.Range("K11").Value = 1093
$ws.Range("M11").Value = -953
$ws.Range("H12").Value = 269
$ws.Range("I12").Value = 327
$ws.Range("K12").Value = 327
$ws.Range("M12").Value = -157
$ws.Range("H28").Value = 720.3570999999999
$ws.Range("I28").Value = 416.81818
$ws.Range("J28").Value = 1833.3334
$ws.Range("K28").Value = 416.81818
$ws.Range("L28").Value = 1833.3334
$ws.Range("M28").Value = 68.18182000000002
$ws.Range("N28").Value = -2803.3334
$ws.Range("H38").Value = 47.285713
$ws.Range("I38").Value = 47.285713
$ws.Range("K38").Value = 141.857139
$ws.Range("M38").Value = 230.142861
$ws.Range("H39").Value = 147.1875
$ws.Range("I39").Value = 147.1875
$ws.Range("K39").Value = 441.5625
$ws.Range("M39").Value = -145.5625
$ws.Range("H55").Value = 523.375
$ws.Range("J55").Value = 607.4
$ws.Range("L55").Value = 607.4
$ws.Range("N55").Value = -1035.4
$ws.Range("H58").Value = 1559.1
$ws.Range("I58").Value = 975
$ws.Range("J58").Value = 1705.125
$ws.Range("K58").Value = 2925
$ws.Range("L58").Value = 5115.375
$ws.Range("M58").Value = -2775
$ws.Range("N58").Value = -5415.375
$ws.Range("H62").Value = 4864.25
$ws.Range("I62").Value = 4025
$ws.Range("J62").Value = 6039.2
$ws.Range("K62").Value = 4025
$ws.Range("L62").Value = 6039.2
$ws.Range("M62").Value = -3401
$ws.Range("N62").Value = -7287.2
$ws.Range("H65").Value = 4864.25
$ws.Range("I65").Value = 4025
$ws.Range("J65").Value = 6039.2
$ws.Range("K65").Value = 20125
$ws.Range("L65").Value = 30196
$ws.Range("M65").Value = -17005
$ws.Range("N65").Value = -36436
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 658.13336
$ws.Range("J88").Value = 664.6667
$ws.Range("L88").Value = 664.6667
$ws.Range("N88").Value = -1476.6667
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 658.13336
$ws.Range("J91").Value = 664.6667
$ws.Range("L91").Value = 664.6667
$ws.Range("N91").Value = -3472.6667
$ws.Range("H107").Value = 1086.75
$ws.Range("I107").Value = 1086.75
$ws.Range("K107").Value = 1086.75
$ws.Range("M107").Value = 833.25
$ws.Range("H111").Value = 965.5454999999999
$ws.Range("I111").Value = 762.1
$ws.Range("K111").Value = 2286.3
$ws.Range("M111").Value = 780.6999999999998
$ws.Range("H137").Value = 1346.2693
$ws.Range("I137").Value = 1141.7273
$ws.Range("J137").Value = 2471.25
$ws.Range("K137").Value = 3425.1819
$ws.Range("L137").Value = 7413.75
$ws.Range("M137").Value = -875.1819
$ws.Range("N137").Value = -12513.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 5350
$ws.Range("I12").Value = 700
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 700
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -527
$ws.Range("N12").Value = -10346
$ws.Range("H32").Value = 6167.2573
$ws.Range("I32").Value = 4878.0586
$ws.Range("K32").Value = 4878.0586
$ws.Range("M32").Value = -4591.0586
$ws.Range("H74").Value = 3100.25
$ws.Range("I74").Value = 3109.4546
$ws.Range("K74").Value = 3109.4546
$ws.Range("M74").Value = -2235.4546
$ws.Range("H77").Value = 3100.25
$ws.Range("I77").Value = 3109.4546
$ws.Range("K77").Value = 15547.273
$ws.Range("M77").Value = -11179.273
$ws.Range("H110").Value = 2455.25
$ws.Range("I110").Value = 1607.3334
$ws.Range("K110").Value = 1607.3334
$ws.Range("M110").Value = 437.6666
$ws.Range("H132").Value = 3457.2856
$ws.Range("I132").Value = 3457.2856
$ws.Range("K132").Value = 10371.8568
$ws.Range("M132").Value = -7841.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 198.66667
$ws.Range("H134").Value = 9000
$ws.Range("I134").Value = 9000
$ws.Range("K134").Value = 27000
$ws.Range("M134").Value = -24465
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2434.3333
$ws.Range("I94").Value = 2439.5
$ws.Range("K94").Value = 2439.5
$ws.Range("M94").Value = -1988.5
$ws.Range("H107").Value = 1832
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2400
$ws.Range("I4").Value = 2362.5
$ws.Range("K4").Value = 7087.5
$ws.Range("M4").Value = -6975.5
$ws.Range("H22").Value = 3815.3333
$ws.Range("J22").Value = 3815.3333
$ws.Range("L22").Value = 11445.9999
$ws.Range("N22").Value = -11783.9999
$ws.Range("H27").Value = 3815.3333
$ws.Range("J27").Value = 3815.3333
$ws.Range("L27").Value = 11445.9999
$ws.Range("N27").Value = -11649.9999
$ws.Range("H44").Value = 248.21428
$ws.Range("J44").Value = 248.21428
$ws.Range("L44").Value = 744.64284
$ws.Range("N44").Value = -1540.64284
$ws.Range("H80").Value = 1987.5
$ws.Range("I80").Value = 1987.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5962.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -5026.5
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1987.5
$ws.Range("I83").Value = 1987.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 17887.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -13207.5
$ws.Range("N83").ClearContents()
$ws.Range("H114").Value = 236.25
$ws.Range("I114").Value = 215
$ws.Range("K114").Value = 645
$ws.Range("M114").Value = 2609
$ws.Range("H117").Value = 969.5
$ws.Range("I117").Value = 492.25
$ws.Range("J117").Value = 1924
$ws.Range("K117").Value = 1476.75
$ws.Range("L117").Value = 5772
$ws.Range("M117").Value = 1965.25
$ws.Range("N117").Value = -12656
$ws.Range("H131").Value = 2281.3572
$ws.Range("J131").Value = 2554.4
$ws.Range("L131").Value = 7663.200000000001
$ws.Range("N131").Value = -17743.2
$ws.Range("H140").Value = 1268.5555
$ws.Range("I140").Value = 852.125
$ws.Range("K140").Value = 2556.375
$ws.Range("M140").Value = 2623.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 50579.777
$ws.Range("J20").Value = 50579.777
$ws.Range("L20").Value = 50579.777
$ws.Range("N20").Value = -51069.777
$ws.Range("H24").Value = 2360857.2
$ws.Range("I24").Value = 35000000
$ws.Range("J24").Value = 29489.785
$ws.Range("K24").Value = 35000000
$ws.Range("L24").Value = 29489.785
$ws.Range("M24").Value = -34999827
$ws.Range("N24").Value = -29835.785
$ws.Range("H113").Value = 1945.7
$ws.Range("I113").Value = 1937.125
$ws.Range("J113").Value = 1980
$ws.Range("K113").Value = 1937.125
$ws.Range("L113").Value = 1980
$ws.Range("M113").Value = 232.875
$ws.Range("N113").Value = -6320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1024.75
$ws.Range("I7").Value = 999.6667
$ws.Range("J7").Value = 1100
$ws.Range("K7").Value = 999.6667
$ws.Range("L7").Value = 1100
$ws.Range("M7").Value = -887.6667
$ws.Range("N7").Value = -1324
$ws.Range("H40").Value = 6639.5
$ws.Range("I40").Value = 6266.222
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 6266.222
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -6130.222
$ws.Range("N40").Value = -10271
$ws.Range("H46").Value = 3582.1765
$ws.Range("I46").Value = 2666.5
$ws.Range("K46").Value = 2666.5
$ws.Range("M46").Value = -2478.5
$ws.Range("H55").Value = 792
$ws.Range("I55").Value = 954.6667
$ws.Range("K55").Value = 954.6667
$ws.Range("M55").Value = -781.6667
$ws.Range("H82").Value = 1484.9
$ws.Range("I82").Value = 1249.8
$ws.Range("J82").Value = 1720
$ws.Range("K82").Value = 1249.8
$ws.Range("L82").Value = 1720
$ws.Range("M82").Value = -888.8
$ws.Range("N82").Value = -2442
$ws.Range("H85").Value = 1484.9
$ws.Range("I85").Value = 1249.8
$ws.Range("J85").Value = 1720
$ws.Range("K85").Value = 1249.8
$ws.Range("L85").Value = 1720
$ws.Range("M85").Value = -1.799999999999955
$ws.Range("N85").Value = -4216
$ws.Range("H122").Value = 3034.3333
$ws.Range("I122").Value = 3108
$ws.Range("J122").Value = 2666
$ws.Range("K122").Value = 9324
$ws.Range("L122").Value = 7998
$ws.Range("M122").Value = -6874
$ws.Range("N122").Value = -12898
$ws.Range("H126").Value = 1024.75
$ws.Range("I126").Value = 999.6667
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 2999.0001
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -529.0001000000002
$ws.Range("N126").Value = -8240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2000
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2226
$ws.Range("H81").Value = 438.33334
$ws.Range("I81").Value = 379.2857
$ws.Range("K81").Value = 758.5714
$ws.Range("M81").Value = 302.4286
$ws.Range("H84").Value = 438.33334
$ws.Range("I84").Value = 379.2857
$ws.Range("K84").Value = 3792.857
$ws.Range("M84").Value = 1511.143
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 373
$ws.Range("N96").ClearContents()
$ws.Range("H117").Value = 24704.5
$ws.Range("J117").Value = 24704.5
$ws.Range("L117").Value = 24704.5
$ws.Range("N117").Value = -33882.5
$ws.Range("H132").Value = 1812.4286
$ws.Range("I132").Value = 1281.1666
$ws.Range("K132").Value = 3843.4998
$ws.Range("M132").Value = -1313.4998

